$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns (N, O), pushing the old "Notes" column from O to P.
$ws.Columns.Item(14).Insert()

# New header row cells for the inserted columns.
$ws.Range("N1").Value = "new_data_dimension"
$ws.Range("O1").Value = "new_memory_GB"

# Row 2 (MedSea): L2/M2 were placeholder text; replace with real numbers/formulas,
# and fill in the new N2/O2 formulas.
$ws.Range("L2").Value = 0
$ws.Range("M2").Formula = "=K2/1000-L2"
$ws.Range("N2").Formula = "=ROUND(1000*M2/F2, 0)*D2"
$ws.Range("O2").Formula = "=(32/8)*N2/(2^30)"

# Rows 3-5 only need the new N/O formulas.
$ws.Range("N3").Formula = "=ROUND(1000*M3/F3, 0)*D3"
$ws.Range("O3").Formula = "=(32/8)*N3/(2^30)"

$ws.Range("N4").Formula = "=ROUND(1000*M4/F4, 0)*D4"
$ws.Range("O4").Formula = "=(32/8)*N4/(2^30)"

$ws.Range("N5").Formula = "=ROUND(1000*M5/F5, 0)*D5"
$ws.Range("O5").Formula = "=(32/8)*N5/(2^30)"

# Row 6 (OOI_Silixa_South): L6/M6 get real numbers/formulas, plus new N6/O6.
$ws.Range("L6").Value = 0
$ws.Range("M6").Formula = "=K6/1000-L6"
$ws.Range("N6").Formula = "=ROUND(1000*M6/F6, 0)*D6"
$ws.Range("O6").Formula = "=(32/8)*N6/(2^30)"

# Row 7 (Svalbard): L7/M7 get real numbers, plus new N7/O7.
$ws.Range("L7").Value = 40
$ws.Range("M7").Value = 100
$ws.Range("N7").Formula = "=ROUND(1000*M7/F7, 0)*D7"
$ws.Range("O7").Formula = "=(32/8)*N7/(2^30)"

# Remove the stray scratch formula in G14.
$ws.Range("G14").Clear()

# Match the saved selection/view from the authored workbook.
$ws.Range("A7").Select()
